$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# New paper rows (20-25), data taken from the commit diff.
# Columns: A=Name B=Index C=Link D=Author E=Year F=Machine Learning
#          G=Mechanism H=Agent Report I=Learner Objective
#          J=Strategy Proofness K=Payment L=Dynamic Game
# ---------------------------------------------------------------------

$rows = @(
    @{ Row=20; Height=57;    A="Causal Strategic Linear Regression";               B=19; C="https://proceedings.icml.cc/static/paper_files/icml/2020/5468-Paper.pdf";                D="Shavit, Edelman, Axelrod";              E=2020; F="Linear regression";  G="NA";                    H="Private x"; I="Loss mimization";    J="No";             K="No";  L="No"  },
    @{ Row=21; Height=42.75; A="Classifiers Induce Invest Effort Strategically";   B=20; C="https://arxiv.org/pdf/1807.05307.pdf";                                                                    D="Kleinberg, Raghavan";                   E=2019; F="Classification";     G="Linear";                H="Hidden Effort"; I="Profit maximization"; J="Not applicable"; K="Yes"; L="No"  },
    @{ Row=22; Height=57;    A="Strategic Noise in Linear Regression";             B=21; C="http://www.cs.toronto.edu/~nisarg/papers/equilibria_linreg.pdf";                                           D="Hossain, Shah";                         E=2020; F="Linear regression";  G="Lp Norm Minimization";  H="Private y"; I="Loss mimization";    J="Yes";            K="No";  L="No"  },
    @{ Row=23; Height=42.75; A="Learning Strategy-Aware Linear Classifiers";       B=22; C="https://arxiv.org/pdf/1911.04004.pdf";                                                                      D="Chen, Liu, Podimata";                   E=2020; F="Linear regression";  G="NA";                    H="Private x"; I="Regret minimization"; J="Not applicable"; K="No";  L="Yes" },
    @{ Row=24; Height=71.25; A="Incentive-Aware PAC Learning";                     B=23; C="https://users.cs.duke.edu/~hrzhang/papers/incentive-aware_learning.pdf";                                   D="Zhang, Conitzer";                       E=2020; F="General";            G="?";                     H="Private x"; I="Loss mimization";    J="Not applicable"; K="No";  L="No"  },
    @{ Row=25; Height=71.25; A="Classification with Strategically Withheld Data";  B=24; C="https://gradanovic.github.io/incentives_in_ML_icml2020_ws/papers/IML2020_paper_13.pdf";                   D="Krishnaswamy, Li, Rein, Zhang, Conitzer"; E=2020; F="Classification";   G="Mincut, Hill-Climbing"; H="Private x"; I="Loss mimization";    J="Yes";            K="No";  L="No"  }
)

foreach ($r in $rows) {
    $rowIdx = $r.Row

    # Bring over the formatting (alignment/wrap + hyperlink style) of the
    # previous last data row (19) onto the whole new row first.
    $ws.Range("A19:L19").Copy()
    $ws.Range("A" + $rowIdx + ":L" + $rowIdx).PasteSpecial(-4122)

    $ws.Cells.Item($rowIdx, 1).Value = $r.A
    $ws.Cells.Item($rowIdx, 2).Value = $r.B
    $ws.Cells.Item($rowIdx, 3).Value = $r.C
    $ws.Cells.Item($rowIdx, 4).Value = $r.D
    $ws.Cells.Item($rowIdx, 5).Value = $r.E
    $ws.Cells.Item($rowIdx, 6).Value = $r.F
    $ws.Cells.Item($rowIdx, 7).Value = $r.G
    $ws.Cells.Item($rowIdx, 8).Value = $r.H
    $ws.Cells.Item($rowIdx, 9).Value = $r.I
    $ws.Cells.Item($rowIdx, 10).Value = $r.J
    $ws.Cells.Item($rowIdx, 11).Value = $r.K
    $ws.Cells.Item($rowIdx, 12).Value = $r.L

    # Register the real hyperlink relationship for column C, then restore
    # the shared Hyperlink cell-style (Add() redefines it to a fresh xf).
    $ws.Hyperlinks.Add($ws.Cells.Item($rowIdx, 3), $r.C) | Out-Null
    $ws.Range("C19").Copy()
    $ws.Cells.Item($rowIdx, 3).PasteSpecial(-4122)

    $ws.Rows.Item($rowIdx).RowHeight = $r.Height
}

# Match the author's final cursor position after entering the new rows.
$ws.Range("A26").Select() | Out-Null
